$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels: the "Question" column header becomes "question" (lowercase)
# and the "Serial_Number" column header becomes "qstNumber".
# (Order matters for how the shared-strings table is rebuilt, so set B1 first.)
$ws.Range("B1").Value = "question"
$ws.Range("A1").Value = "qstNumber"

# Move the active selection to B9 (matches updated <selection> in sheet XML).
$ws.Range("B9").Select()
